# Natmi following Dr Hou advice
# Recomputed LR-pair stats for Fgf2-Gpc4: the "Sending cluster" set grows
# from {FAPs, M1, sCs} to {ECs, FAPs, M1, sCs} (4 senders x 3 targets = 12
# rows instead of 3 senders x 3 targets = 9), and every numeric column is
# refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the 16 numeric statistic columns.
$rows = @(
    @("ECs",  "Fgf2", "Gpc4", "ECs",  2, 0.6666666666666666, 0.8896923333333334, 2.669077,           0.09326752961394506, 0.09326752961394505, 3, 1, 15.740393,          47.221179,          0.2847568403735705, 0.2847568403735705, 14.00410697575367, 126.036962781783,  0.02655856704231541,  0.02655856704231541),
    @("ECs",  "Fgf2", "Gpc4", "FAPs", 2, 0.6666666666666666, 0.8896923333333334, 2.669077,           0.09326752961394506, 0.09326752961394505, 3, 1, 33.82224299999999, 101.466729,         0.6118725911752718, 0.6118725911752717, 30.091390293237,   270.822512639133,  0.05706784501740097,  0.05706784501740095),
    @("ECs",  "Fgf2", "Gpc4", "sCs",  2, 0.6666666666666666, 0.8896923333333334, 2.669077,           0.09326752961394506, 0.09326752961394505, 3, 1, 5.713974666666666, 17.141924,          0.1033705684511578, 0.1033705684511578, 5.083679453794223, 45.753115084148,   0.009641117554228699, 0.009641117554228695),
    @("FAPs", "Fgf2", "Gpc4", "ECs",  3, 1,                  7.418580000000001, 22.25574,           0.777698766101638,   0.777698766101638,   3, 1, 15.740393,          47.221179,          0.2847568403735705, 0.2847568403735705, 116.77136470194,   1050.94228231746,  0.2214550433975269,   0.2214550433975269),
    @("FAPs", "Fgf2", "Gpc4", "FAPs", 3, 1,                  7.418580000000001, 22.25574,           0.777698766101638,   0.777698766101638,   3, 1, 33.82224299999999, 101.466729,         0.6118725911752718, 0.6118725911752717, 250.91301547494,   2258.21713927446,  0.4758525591684209,   0.4758525591684208),
    @("FAPs", "Fgf2", "Gpc4", "sCs",  3, 1,                  7.418580000000001, 22.25574,           0.777698766101638,   0.777698766101638,   3, 1, 5.713974666666666, 17.141924,          0.1033705684511578, 0.1033705684511578, 42.38957818264,    381.50620364376,   0.08039116353569035,  0.08039116353569034),
    @("M1",   "Fgf2", "Gpc4", "ECs",  3, 1,                  0.298413,          0.895239,           0.031282997809377,   0.03128299780937701, 3, 1, 15.740393,          47.221179,          0.2847568403735705, 0.2847568403735705, 4.697137896308999, 42.274241066781,   0.008908047613611522, 0.008908047613611523),
    @("M1",   "Fgf2", "Gpc4", "FAPs", 3, 1,                  0.298413,          0.895239,           0.031282997809377,   0.03128299780937701, 3, 1, 33.82224299999999, 101.466729,         0.6118725911752718, 0.6118725911752717, 10.092997000359,   90.83697300323099, 0.01914120892935386,  0.01914120892935386),
    @("M1",   "Fgf2", "Gpc4", "sCs",  3, 1,                  0.298413,          0.895239,           0.031282997809377,   0.03128299780937701, 3, 1, 5.713974666666666, 17.141924,          0.1033705684511578, 0.1033705684511578, 1.705124322204,    15.346118899836,   0.003233741266411626, 0.003233741266411626),
    @("sCs",  "Fgf2", "Gpc4", "ECs",  3, 1,                  0.9324580000000001,2.797374,           0.09775070647503986, 0.09775070647503986, 3, 1, 15.740393,          47.221179,          0.2847568403735705, 0.2847568403735705, 14.677255375994,   132.095298383946,  0.02783518232011667,  0.02783518232011667),
    @("sCs",  "Fgf2", "Gpc4", "FAPs", 3, 1,                  0.9324580000000001,2.797374,           0.09775070647503986, 0.09775070647503986, 3, 1, 33.82224299999999, 101.466729,         0.6118725911752718, 0.6118725911752717, 31.537821063294,   283.840389569646,  0.05981097806009606,  0.05981097806009605),
    @("sCs",  "Fgf2", "Gpc4", "sCs",  3, 1,                  0.9324580000000001,2.797374,           0.09775070647503986, 0.09775070647503986, 3, 1, 5.713974666666666, 17.141924,          0.1033705684511578, 0.1033705684511578, 5.328041389730667, 47.95237250757601, 0.01010454609482715,  0.01010454609482714)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
